$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a paragraph index by a distinctive substring of its text.
# ---------------------------------------------------------------------------
function Find-ParagraphIndexByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: insert a brand-new bullet item, at the very top of the
# "Puesta en funcionamiento" numbered list, instructing the reader to clone
# the repository. The paragraph that used to be first in that list
# ("Situate en el directorio raiz del proyecto...") is pushed down and
# keeps all of its existing text/formatting untouched.
# ---------------------------------------------------------------------------
$firstListItemIdx = Find-ParagraphIndexByText $d "Sitúate en el directorio raíz del proyecto"
$firstListItem = $d.Paragraphs.Item($firstListItemIdx)

# Inserting a paragraph break before the first item duplicates its pPr
# (style + numbering), giving the new paragraph the same bullet formatting.
$firstListItem.Range.InsertParagraphBefore()

$newFirstItem = $d.Paragraphs.Item($firstListItemIdx)
$newFirstItem.Range.Text = "Haz un clone del proyecto situado en mi repositorio: https://github.com/coviza/ProyectoFinal"

# ---------------------------------------------------------------------------
# Change 2: add a new, empty "Parrafo de lista" paragraph right after the
# last bullet of that same list ("Podras acceder desde..."), before the
# pre-existing blank separator paragraph that leads into the "Links:" block.
# ---------------------------------------------------------------------------
$lastListItemIdx = Find-ParagraphIndexByText $d "Podrás acceder desde"
$lastListItem = $d.Paragraphs.Item($lastListItemIdx)

# Find the blank separator paragraph that immediately follows the list so we
# can insert the new paragraph right before it (i.e. right after the list).
$separatorIdx = $lastListItemIdx + 1
$separator = $d.Paragraphs.Item($separatorIdx)
$separator.Range.InsertParagraphBefore()

$newLastParagraph = $d.Paragraphs.Item($separatorIdx)
# Apply the "List Paragraph" style only, without any bullet/numbering, which
# matches assigning the "Prrafodelista" paragraph style with no numPr.
$newLastParagraph.Style = "Prrafodelista"
